$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ranges that must stay plain text (so numeric-looking / date-looking
# strings are not auto-converted by Excel into numbers or dates).
$colA = $ws.Range("A2:A3")
$colEF = $ws.Range("E2:F3")

$colA.NumberFormat = "@"
$colEF.NumberFormat = "@"

$ws.Range("A2").Value = "6787743305"
$ws.Range("B2").Value = "Ascend Laboratories, LLC"
$ws.Range("C2").Value = "Aripiprazole"
$ws.Range("D2").Value = "15 mg/1"
$ws.Range("E2").Value = "22140477"
$ws.Range("F2").Value = "01/24/31"
$ws.Range("G2").Value = 6

$ws.Range("A3").Value = "6787743305"
$ws.Range("B3").Value = "Ascend Laboratories, LLC"
$ws.Range("C3").Value = "Aripiprazole"
$ws.Range("D3").Value = "15 mg/1"
$ws.Range("E3").Value = "22140477"
$ws.Range("F3").Value = "01/24/31"
$ws.Range("G3").Value = 99

# Restore the default style so no extra style/number-format definitions
# linger on the cells while keeping their values stored as text.
$colA.Style = "Normal"
$colEF.Style = "Normal"
